$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'50.999.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  -1.21%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Formula = "'3.043.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +0.74%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Formula = "'  +0.16%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Formula = "'385.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  +1.65%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Formula = "'101.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  -0.88%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Formula = "'0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  -2.41%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Formula = "'  +0.01%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Formula = "'0.575"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  -2.44%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Formula = "'36.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  -0.89%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Formula = "'  -0.11%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Formula = "'0.0844"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  -2.20%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Formula = "'3.537.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  +1.09%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Formula = "'18.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  -1.27%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Formula = "'7.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  -1.24%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Formula = "'3.108.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +2.72%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Formula = "'0.979"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +0.59%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Formula = "'10.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +0.25%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Formula = "'51.045.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  -1.08%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Formula = "'3.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +2.64%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Formula = "'0.0₃0949"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -1.54%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Formula = "'12.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  -2.32%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Formula = "'69.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -0.83%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Formula = "'262.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  -1.59%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Formula = "'3.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  -1.57%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Formula = "'7.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  -5.56%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Formula = "'26.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  +2.19%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Formula = "'Dai"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Formula = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Formula = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  -0.07%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Formula = "'RenderToken"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Formula = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Formula = "'7.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  -5.10%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Formula = "'0.162"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  -6.01%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Formula = "'  -3.76%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Formula = "'10.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  +0.84%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Formula = "'35.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  +3.51%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Formula = "'0.0465"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  +3.61%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Formula = "'  -2.82%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Formula = "'50.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  -1.04%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Formula = "'  +0.09%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Formula = "'  +1.91%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Formula = "'  -1.88%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Formula = "'129.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +1.54%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Formula = "'  -1.92%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Formula = "'  -1.39%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Formula = "'3.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  -0.12%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Formula = "'16.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  -3.86%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Formula = "'2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -3.41%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Formula = "'21.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "'  -0.23%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Formula = "'2.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  +2.79%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Formula = "'2.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  -0.15%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Formula = "'2.049.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  +1.14%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Formula = "'9.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +9.02%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Formula = "'BEAM"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Formula = "'https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Formula = "'0.0319"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +0.63%  "
$ws.Range("E51").Style = "Normal"
